$d = $word.ActiveDocument

# 1. Expand the "5mn-30mn" sentence into the longer revised text.
$old1 = "Write a script for the favorite coins that show a price movement OPPOSITE to bitcoin within a 5mn-30mn range. (Then Jim you need to keep an eye on them…)"
$new1 = "Write a script for the favorite coins that show a price movement OPPOSITE to bitcoin within a 5mn-60mn range. Also check for bullish or bearish movement of those coins when bitcoin is consolidating in the within 5-60mn in the 5mn time frame (you need to write a very clear elaborate script). You have to have an automated sensitive reliable system to catch that anomaly. "
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

Write-Output "done"
